{"js": "// Office.js (Word JavaScript API) script\n// Applies the LOM3099 \"Est\u00e1tica\" syllabus updates described by the diff.\n\nconst replacements = [\n  {\n    search: \"Ativa\u00e7\u00e3o: 01/01/2020\",\n    replace: \"Ativa\u00e7\u00e3o: 01/01/2022\"\n  },\n  {\n    search: \"Est\u00e1tica de Part\u00edculas. Est\u00e1tica de Corpos R\u00edgidos.An\u00e1lise de Estruturas.\",\n    replace: \"Est\u00e1tica de Part\u00edculas. Est\u00e1tica de Corpos R\u00edgidos. Equil\u00edbrio de Corpos R\u00edgidos. An\u00e1lise de Estruturas.\"\n  },\n  {\n    search: \"Para compor a Nota no Semestre (NS) ser\u00e3o feitas duas avalia\u00e7\u00f5es (P1 e P2).\",\n    replace: \"Os alunos ser\u00e3o avaliados continuamente a qual ser\u00e3o considerados: provas escritas, exerc\u00edcios propostos e semin\u00e1rios\"\n  },\n  {\n    search: \"NS = NP1 + NP2. NP1: quest\u00f5es da P1 valendo at\u00e9 4p. no total; NP2: quest\u00f5es da P2 valendo at\u00e9 6p. no total. Ser\u00e3o considerados aprovados os alunos que obtiverem: NS maior ou igual a 5,0. Ser\u00e3o considerados reprovados os alunos que obtiverem: NS menor que 3,0. Para os alunos em que NS \u00e9 maior ou igual a 3,0 e menor que 5,0 ser\u00e1 dada uma prova de recupera\u00e7\u00e3o (R).\",\n    replace: \"Para compor a Nota Final (NF) ser\u00e3o consideradas as avalia\u00e7\u00f5es de Provas Escritas (P1 e P2) e Exerc\u00edcios Propostos e Semin\u00e1rios (T) em que:NF = (P1 + P2 + T)/3.  Ser\u00e3o considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Ser\u00e3o considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS \u00e9 maior ou igual a 3,0 e menor que 5,0 ser\u00e1 dada uma prova de recupera\u00e7\u00e3o (R).\"\n  },\n  {\n    search: \"A prova de Recupera\u00e7\u00e3o (R) ir\u00e1 compor a nota final (NF) da seguinte forma: NF = (R + NS)/2. Ser\u00e3o considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0.\",\n    replace: \": A prova de Recupera\u00e7\u00e3o (R) ir\u00e1 compor a nota final de recupera\u00e7\u00e3o (NR) da seguinte forma: NR = (R + NF)/2. Ser\u00e3o considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0.\"\n  },\n  {\n    search: \"1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Est\u00e1tica e Mec\u00e2nica dos Materiais. S\u00e3o Paulo: McGraw Hill, 2013, 728p. 2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mec\u00e2nica vetorial para engenheiros: est\u00e1tica./ S\u00e3o Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mec\u00e2nica para engenharia vol.1: est\u00e1tica. S\u00e3o Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mec\u00e2nica para engenharia \u2013 Est\u00e1tica. Grupo GEN Editora LTC, 6a Ed., 2009, 364p.5. RUIZ, C.C.de La P. Fundamentos de mec\u00e2nica para engenharia \u2013 Est\u00e1tica. Grupo GEN Editora LTC, 2017, 306p.\",\n    replace: \"1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Est\u00e1tica e Mec\u00e2nica dos Materiais. S\u00e3o Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mec\u00e2nica vetorial para engenheiros: Est\u00e1tica. S\u00e3o Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mec\u00e2nica para engenharia vol.1: est\u00e1tica. S\u00e3o Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mec\u00e2nica para engenharia \u2013 Est\u00e1tica. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mec\u00e2nica para engenharia \u2013 Est\u00e1tica. Grupo GEN Editora LTC, 2017, 306p.\"\n  }\n];\n\nfor (const { search, replace } of replacements) {\n  const results = context.document.body.search(search, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop script\n# Applies the LOM3099 \"Est\u00e1tica\" syllabus updates described by the diff.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1  # wdFindContinue\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2020\" \"Ativa\u00e7\u00e3o: 01/01/2022\"\n\nReplace-Text \"Est\u00e1tica de Part\u00edculas. Est\u00e1tica de Corpos R\u00edgidos.An\u00e1lise de Estruturas.\" \"Est\u00e1tica de Part\u00edculas. Est\u00e1tica de Corpos R\u00edgidos. Equil\u00edbrio de Corpos R\u00edgidos. An\u00e1lise de Estruturas.\"\n\nReplace-Text \"Para compor a Nota no Semestre (NS) ser\u00e3o feitas duas avalia\u00e7\u00f5es (P1 e P2).\" \"Os alunos ser\u00e3o avaliados continuamente a qual ser\u00e3o considerados: provas escritas, exerc\u00edcios propostos e semin\u00e1rios\"\n\nReplace-Text \"NS = NP1 + NP2. NP1: quest\u00f5es da P1 valendo at\u00e9 4p. no total; NP2: quest\u00f5es da P2 valendo at\u00e9 6p. no total. Ser\u00e3o considerados aprovados os alunos que obtiverem: NS maior ou igual a 5,0. Ser\u00e3o considerados reprovados os alunos que obtiverem: NS menor que 3,0. Para os alunos em que NS \u00e9 maior ou igual a 3,0 e menor que 5,0 ser\u00e1 dada uma prova de recupera\u00e7\u00e3o (R).\" \"Para compor a Nota Final (NF) ser\u00e3o consideradas as avalia\u00e7\u00f5es de Provas Escritas (P1 e P2) e Exerc\u00edcios Propostos e Semin\u00e1rios (T) em que:NF = (P1 + P2 + T)/3.  Ser\u00e3o considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Ser\u00e3o considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS \u00e9 maior ou igual a 3,0 e menor que 5,0 ser\u00e1 dada uma prova de recupera\u00e7\u00e3o (R).\"\n\nReplace-Text \"A prova de Recupera\u00e7\u00e3o (R) ir\u00e1 compor a nota final (NF) da seguinte forma: NF = (R + NS)/2. Ser\u00e3o considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0.\" \": A prova de Recupera\u00e7\u00e3o (R) ir\u00e1 compor a nota final de recupera\u00e7\u00e3o (NR) da seguinte forma: NR = (R + NF)/2. Ser\u00e3o considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0.\"\n\nReplace-Text \"1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Est\u00e1tica e Mec\u00e2nica dos Materiais. S\u00e3o Paulo: McGraw Hill, 2013, 728p. 2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mec\u00e2nica vetorial para engenheiros: est\u00e1tica./ S\u00e3o Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mec\u00e2nica para engenharia vol.1: est\u00e1tica. S\u00e3o Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mec\u00e2nica para engenharia \u2013 Est\u00e1tica. Grupo GEN Editora LTC, 6a Ed., 2009, 364p.5. RUIZ, C.C.de La P. Fundamentos de mec\u00e2nica para engenharia \u2013 Est\u00e1tica. Grupo GEN Editora LTC, 2017, 306p.\" \"1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Est\u00e1tica e Mec\u00e2nica dos Materiais. S\u00e3o Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mec\u00e2nica vetorial para engenheiros: Est\u00e1tica. S\u00e3o Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mec\u00e2nica para engenharia vol.1: est\u00e1tica. S\u00e3o Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mec\u00e2nica para engenharia \u2013 Est\u00e1tica. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mec\u00e2nica para engenharia \u2013 Est\u00e1tica. Grupo GEN Editora LTC, 2017, 306p.\"\n"}
